$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 868.3889
$ws.Range("J19").Value = 1089.6923
$ws.Range("L19").Value = 1089.6923
$ws.Range("N19").Value = -1439.6923

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3850
$ws.Range("J51").Value = 4175
$ws.Range("L51").Value = 4175
$ws.Range("N51").Value = -5143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1385.5897
$ws.Range("I132").Value = 1010.17645
$ws.Range("K132").Value = 3030.52935
$ws.Range("M132").Value = -500.5293500000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5418.7354
$ws.Range("I137").Value = 3322.375
$ws.Range("K137").Value = 9967.125
$ws.Range("M137").Value = -7417.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3710611.5
$ws.Range("I32").Value = 4004960.5
$ws.Range("K32").Value = 4004960.5
$ws.Range("M32").Value = -4004673.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5328.357
$ws.Range("J45").Value = 5980.5454
$ws.Range("L45").Value = 5980.5454
$ws.Range("N45").Value = -6734.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 6000
$ws.Range("I57").Value = 6000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 50004044
$ws.Range("I61").Value = 2408.3333
$ws.Range("J61").Value = 125006500
$ws.Range("K61").Value = 2408.3333
$ws.Range("L61").Value = 125006500
$ws.Range("M61").Value = -2196.3333
$ws.Range("N61").Value = -125006924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 31348.97
$ws.Range("I74").Value = 43577.082
$ws.Range("J74").Value = 4669.4546
$ws.Range("K74").Value = 43577.082
$ws.Range("L74").Value = 4669.4546
$ws.Range("M74").Value = -42703.082
$ws.Range("N74").Value = -6417.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 31348.97
$ws.Range("I77").Value = 43577.082
$ws.Range("J77").Value = 4669.4546
$ws.Range("K77").Value = 217885.41
$ws.Range("L77").Value = 23347.273
$ws.Range("M77").Value = -213517.41
$ws.Range("N77").Value = -32083.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1726.1111
$ws.Range("J88").Value = 1816.875
$ws.Range("L88").Value = 1816.875
$ws.Range("N88").Value = -2628.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1726.1111
$ws.Range("J91").Value = 1816.875
$ws.Range("L91").Value = 1816.875
$ws.Range("N91").Value = -4624.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 44727.4
$ws.Range("J96").Value = 44727.4
$ws.Range("L96").Value = 44727.4
$ws.Range("N96").Value = -50219.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 50004044
$ws.Range("I136").Value = 2408.3333
$ws.Range("J136").Value = 125006500
$ws.Range("K136").Value = 7224.999899999999
$ws.Range("L136").Value = 375019500
$ws.Range("M136").Value = -4674.999899999999
$ws.Range("N136").Value = -375024600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3149.1538
$ws.Range("I105").Value = 2221.5417
$ws.Range("K105").Value = 2221.5417
$ws.Range("M105").Value = -474.5417000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5260.5
$ws.Range("J29").Value = 5260.5
$ws.Range("L29").Value = 5260.5
$ws.Range("N29").Value = -5846.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2031
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 11982.5
$ws.Range("I35").Value = 1125
$ws.Range("K35").Value = 1125
$ws.Range("M35").Value = -831

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 49722.5
$ws.Range("J36").Value = 49722.5
$ws.Range("L36").Value = 49722.5
$ws.Range("N36").Value = -50498.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 50000
$ws.Range("I39").Value = 50000
$ws.Range("K39").Value = 50000
$ws.Range("M39").Value = -49609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 49722.5
$ws.Range("J40").Value = 49722.5
$ws.Range("L40").Value = 49722.5
$ws.Range("N40").Value = -50042.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 44915.668
$ws.Range("J44").Value = 44915.668
$ws.Range("L44").Value = 44915.668
$ws.Range("N44").Value = -45799.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 1000000
$ws.Range("J47").Value = 1000000
$ws.Range("L47").Value = 1000000
$ws.Range("N47").Value = -1001132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 50000
$ws.Range("I49").Value = 50000
$ws.Range("K49").Value = 50000
$ws.Range("M49").Value = -49818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 94036
$ws.Range("J50").Value = 94036
$ws.Range("L50").Value = 94036
$ws.Range("N50").Value = -95286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 31085.363
$ws.Range("J51").Value = 33684.9
$ws.Range("L51").Value = 33684.9
$ws.Range("N51").Value = -35156.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 99776.5
$ws.Range("J59").Value = 99776.5
$ws.Range("L59").Value = 99776.5
$ws.Range("N59").Value = -102066.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 36071.45
$ws.Range("I60").Value = 21916.334
$ws.Range("K60").Value = 21916.334
$ws.Range("M60").Value = -21405.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 31085.363
$ws.Range("J61").Value = 33684.9
$ws.Range("L61").Value = 33684.9
$ws.Range("N61").Value = -34380.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 25000
$ws.Range("K81").Value = 25000
$ws.Range("M81").Value = -24002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 43000
$ws.Range("I82").Value = 43000
$ws.Range("K82").Value = 43000
$ws.Range("M82").Value = -42639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 25000
$ws.Range("K84").Value = 75000
$ws.Range("M84").Value = -70008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 43000
$ws.Range("I85").Value = 43000
$ws.Range("K85").Value = 43000
$ws.Range("M85").Value = -41752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4763160.5
$ws.Range("I105").Value = 7937301
$ws.Range("K105").Value = 7937301
$ws.Range("M105").Value = -7935554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5910.5654
$ws.Range("I132").Value = 5396.4
$ws.Range("J132").Value = 6053.3887
$ws.Range("K132").Value = 16189.2
$ws.Range("L132").Value = 18160.1661
$ws.Range("M132").Value = -13659.2
$ws.Range("N132").Value = -23220.1661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5644.4443
$ws.Range("I134").Value = 4500
$ws.Range("K134").Value = 13500
$ws.Range("M134").Value = -10965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1574.75
$ws.Range("I3").Value = 1574.75
$ws.Range("K3").Value = 4724.25
$ws.Range("M3").Value = -4612.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 96895.81
$ws.Range("I140").Value = 101590.6
$ws.Range("K140").Value = 304771.8
$ws.Range("M140").Value = -299591.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 80000
$ws.Range("J52").Value = 80000
$ws.Range("L52").Value = 80000
$ws.Range("N52").Value = -80518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2227.2856
$ws.Range("I132").Value = 1835.6522
$ws.Range("K132").Value = 5506.9566
$ws.Range("M132").Value = -2976.9566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 26427.46
$ws.Range("J136").Value = 25777.908
$ws.Range("L136").Value = 77333.724
$ws.Range("N136").Value = -82433.724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2759.4
$ws.Range("I16").Value = 2949.25
$ws.Range("K16").Value = 2949.25
$ws.Range("M16").Value = -2779.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("I5").Value = 10000000
$ws.Range("K5").Value = 10000000
$ws.Range("M5").Value = -9999888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5076.3335
$ws.Range("J74").Value = 5076.3335
$ws.Range("L74").Value = 5076.3335
$ws.Range("N74").Value = -6948.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 5076.3335
$ws.Range("J77").Value = 5076.3335
$ws.Range("L77").Value = 15229.0005
$ws.Range("N77").Value = -24589.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 38665
$ws.Range("J82").Value = 38665
$ws.Range("L82").Value = 38665
$ws.Range("N82").Value = -39431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 38665
$ws.Range("J85").Value = 38665
$ws.Range("L85").Value = 38665
$ws.Range("N85").Value = -41317
